# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted above the existing row 54
# ("Vega Modelo de Temuco" / "Arándano (blue)"), pushing the previous
# rows 54-86 down to 55-87 and growing the used range from A1:T86 to
# A1:T87.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 54 (shifts rows 54:86 down to 55:87).
$ws.Rows("54:54").Insert()

# Populate the newly inserted row 54 with the new weekly record.
$ws.Range("A54").Value = 10
$ws.Range("B54").Value = 'Vega Modelo de Temuco'
$ws.Range("C54").Value = 'La Araucanía'
$ws.Range("D54").Value = 44606
$ws.Range("E54").Value = 9
$ws.Range("F54").Value = 'Fruta'
$ws.Range("G54").Value = 100101
$ws.Range("H54").Value = 'Berries'
$ws.Range("I54").Value = 100101001
$ws.Range("J54").Value = 'Arándano (blue)'
$ws.Range("K54").Value = 'Sin especificar'
$ws.Range("L54").Value = 'Primera'
$ws.Range("M54").Value = 300
$ws.Range("N54").Value = 2000
$ws.Range("O54").Value = 2000
$ws.Range("P54").Value = 2000
$ws.Range("Q54").Value = '$/kilo'
$ws.Range("R54").Value = 'Región del Maule'
$ws.Range("S54").Value = 2000
$ws.Range("T54").Value = 1
